# Premier League Project PPT.pptx - "uploading final VS code" edit
#
# 1) Slide 5 (sldId 259): reposition "Picture 4" (the player-stats image) up
#    and to the left.
# 2) Slide 14 (the "What If?" / Q&A slide): the "Open For Questions" run and
#    the following ":" run get merged into a single run reading
#    "Open For Questions:" (keeping the formatting/rPr of the second run).

$p = $ppt.ActivePresentation

# --- 1) Move the picture on slide 5 -----------------------------------
$slide5 = $p.Slides.Item(5)
$pic = $null
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    if ($slide5.Shapes.Item($i).Name -eq "Picture 4") {
        $pic = $slide5.Shapes.Item($i)
    }
}
$pic.Left = 40.33740157480315
$pic.Top = 103.19811023622047

# --- 2) Merge "Open For Questions" + ":" runs on slide 14 --------------
$slide14 = $p.Slides.Item(14)
$shp = $slide14.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# This textbox autofits its height to content (a:spAutoFit); remember the
# current height so it can be restored after the run-merge below (merging
# two runs with identical rendered text into one shouldn't move the box,
# but the engine recomputes the autofit height on any text-content edit).
$origHeight = $shp.Height

$target = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $ptext = $para.Text.TrimEnd([char]13, [char]10)
    if ($ptext -eq "Open For Questions:") {
        $target = $para
    }
}

$paraStart = $target.Start
$firstRun = $tr.Characters($paraStart, 18)
$firstRun.Text = ""

$remainder = $tr.Characters($paraStart, 1)
$remainder.Text = "Open For Questions:"

# Restore the box height (nudge past the Single/float rounding boundary so
# the re-serialized EMU value matches the pre-edit value exactly).
$shp.Height = $origHeight + 0.00002
